# 6.2.1 worksheet update: split the single "2018" data column into three
# year columns (2014 / 2018 / 2023) and refresh the footnote that documents
# the source years.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Make room for two new columns (E, F) to the right of the existing
#    data column D. Excel's column-insert copies the formatting of the
#    column to the left, so E and F pick up D's styles automatically.
# ---------------------------------------------------------------------
$ws.Columns.Item(5).Insert()
$ws.Columns.Item(6).Insert()

# Rows 3 and 31 only ever had a single trailing data cell (D3 / D31) - the
# insert step above spuriously extends them with blank-but-styled E/F
# cells; clear those back out so the row only spans as far as it should.
$ws.Range("E3").Clear()
$ws.Range("F3").Clear()
$ws.Range("E31").Clear()
$ws.Range("F31").Clear()

# ---------------------------------------------------------------------
# 2) The old column D held the 2018 figures; those move to the new
#    column E, then D/F receive the new 2014 / 2023 figures.
# ---------------------------------------------------------------------
$ws.Range("D4:D30").Copy($ws.Range("E4:E30"))

# Year header row
$ws.Range("D4").Value = 2014
$ws.Range("E4").Value = 2018
$ws.Range("F4").Value = 2023

# Data rows: column D (2014) and column F (2023). Column E (2018) already
# holds the figures copied from the old column D above.
$ws.Range("D5").Value = 94.5
$ws.Range("F5").Value = 97.6

$ws.Range("D7").Value = 96.8
$ws.Range("F7").Value = 96.7

$ws.Range("D8").Value = 93
$ws.Range("F8").Value = 98

$ws.Range("D10").Value = 91.8
$ws.Range("F10").Value = 97

$ws.Range("D11").Value = 93
$ws.Range("F11").Value = 98.1

$ws.Range("D12").Value = 95.2
$ws.Range("F12").Value = 98.7

$ws.Range("D13").Value = 97.3
$ws.Range("F13").Value = 99.5

$ws.Range("D14").Value = 94.3
$ws.Range("F14").Value = 97.3

$ws.Range("D15").Value = 91
$ws.Range("F15").Value = 97.1

$ws.Range("D16").Value = 92.7
$ws.Range("F16").Value = 99.1

$ws.Range("D17").Value = 99.1
$ws.Range("F17").Value = 97.8

$ws.Range("D18").Value = 92.4
$ws.Range("F18").Value = 90.4

$ws.Range("D20").Value = 95.1
$ws.Range("F20").Value = 91.1

$ws.Range("D21").Value = 90.6
$ws.Range("F21").Value = 97.5

$ws.Range("D22").Value = 93.6
$ws.Range("F22").Value = 97.6

$ws.Range("D23").Value = 95.6
$ws.Range("F23").Value = 97.9

$ws.Range("D24").Value = 97.4
$ws.Range("F24").Value = 98

$ws.Range("D26").Value = 88.2
$ws.Range("F26").Value = 95.1

$ws.Range("D27").Value = 93.4
$ws.Range("F27").Value = 97.7

$ws.Range("D28").Value = 95.1
$ws.Range("F28").Value = 98.8

$ws.Range("D29").Value = 95.8
$ws.Range("F29").Value = 99

$ws.Range("D30").Value = 97.9
$ws.Range("F30").Value = 97.6

# ---------------------------------------------------------------------
# 3) Refresh the footnote (row 31) so it references all three survey
#    years instead of just 2018.
# ---------------------------------------------------------------------
$ws.Range("A31").Value = " Көп көрсөткүчтүү кластердик изилдөөнүн маалыматтары боюнча, 2014- ж., 2018-ж., 2023-ж."
$ws.Range("B31").Value = "По данным кластерного обследования по многим показателям, 2014г., 2018г., 2023г."
$ws.Range("C31").Value = "According to Multiple Indicator Cluster Survey, 2014, 2018, 2023."

# ---------------------------------------------------------------------
# 4) Cosmetic tweaks that came along with the data change.
# ---------------------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 67.5
$ws.Range("A1").Select()
